# Update NATMI LR-pair metrics (Tnf-Tnfrsf1b) with new TPM-derived values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.624526
$ws.Range("H2").Value = 7.873578
$ws.Range("I2").Value = 0.06442870872595916
$ws.Range("J2").Value = 0.06442870872595916
$ws.Range("M2").Value = 16.05260533333333
$ws.Range("N2").Value = 48.157816
$ws.Range("O2").Value = 0.1752915379534001
$ws.Range("P2").Value = 0.1752915379534001
$ws.Range("Q2").Value = 42.13048006507199
$ws.Range("R2").Value = 379.174320585648
$ws.Range("S2").Value = 0.01129380744092503
$ws.Range("T2").Value = 0.01129380744092503
$ws.Range("G3").Value = 2.624526
$ws.Range("H3").Value = 7.873578
$ws.Range("I3").Value = 0.06442870872595916
$ws.Range("J3").Value = 0.06442870872595916
$ws.Range("O3").Value = 0.07888758308485012
$ws.Range("P3").Value = 0.07888758308485012
$ws.Range("Q3").Value = 18.960252076866
$ws.Range("R3").Value = 170.642268691794
$ws.Range("S3").Value = 0.005082625112668711
$ws.Range("T3").Value = 0.005082625112668711
$ws.Range("G4").Value = 2.624526
$ws.Range("H4").Value = 7.873578
$ws.Range("I4").Value = 0.06442870872595916
$ws.Range("J4").Value = 0.06442870872595916
$ws.Range("M4").Value = 2.098187333333334
$ws.Range("N4").Value = 6.294562000000001
$ws.Range("O4").Value = 0.02291182502385553
$ws.Range("P4").Value = 0.02291182502385553
$ws.Range("Q4").Value = 5.506747209204001
$ws.Range("R4").Value = 49.56072488283601
$ws.Range("S4").Value = 0.00147617930084213
$ws.Range("T4").Value = 0.00147617930084213
$ws.Range("G5").Value = 2.624526
$ws.Range("H5").Value = 7.873578
$ws.Range("I5").Value = 0.06442870872595916
$ws.Range("J5").Value = 0.06442870872595916
$ws.Range("M5").Value = 66.20156266666667
$ws.Range("N5").Value = 198.604688
$ws.Range("O5").Value = 0.7229090539378943
$ws.Range("P5").Value = 0.7229090539378942
$ws.Range("Q5").Value = 173.747722459296
$ws.Range("R5").Value = 1563.729502133664
$ws.Range("S5").Value = 0.04657609687152329
$ws.Range("T5").Value = 0.04657609687152328
$ws.Range("I6").Value = 0.01210207677934112
$ws.Range("J6").Value = 0.01210207677934112
$ws.Range("M6").Value = 16.05260533333333
$ws.Range("N6").Value = 48.157816
$ws.Range("O6").Value = 0.1752915379534001
$ws.Range("P6").Value = 0.1752915379534001
$ws.Range("Q6").Value = 7.913650833305778
$ws.Range("R6").Value = 71.222857499752
$ws.Range("S6").Value = 0.002121391651080836
$ws.Range("T6").Value = 0.002121391651080836
$ws.Range("I7").Value = 0.01210207677934112
$ws.Range("J7").Value = 0.01210207677934112
$ws.Range("O7").Value = 0.07888758308485012
$ws.Range("P7").Value = 0.07888758308485012
$ws.Range("S7").Value = 0.0009547035874295082
$ws.Range("T7").Value = 0.0009547035874295082
$ws.Range("I8").Value = 0.01210207677934112
$ws.Range("J8").Value = 0.01210207677934112
$ws.Range("M8").Value = 2.098187333333334
$ws.Range("N8").Value = 6.294562000000001
$ws.Range("O8").Value = 0.02291182502385553
$ws.Range("P8").Value = 0.02291182502385553
$ws.Range("Q8").Value = 1.034369287357111
$ws.Range("R8").Value = 9.309323586214001
$ws.Range("S8").Value = 0.0002772806655935288
$ws.Range("T8").Value = 0.0002772806655935288
$ws.Range("I9").Value = 0.01210207677934112
$ws.Range("J9").Value = 0.01210207677934112
$ws.Range("M9").Value = 66.20156266666667
$ws.Range("N9").Value = 198.604688
$ws.Range("O9").Value = 0.7229090539378943
$ws.Range("P9").Value = 0.7229090539378942
$ws.Range("Q9").Value = 32.63620083372623
$ws.Range("R9").Value = 293.725807503536
$ws.Range("S9").Value = 0.008748700875237249
$ws.Range("T9").Value = 0.008748700875237249
$ws.Range("I10").Value = 0.003429134645952472
$ws.Range("J10").Value = 0.003429134645952472
$ws.Range("M10").Value = 16.05260533333333
$ws.Range("N10").Value = 48.157816
$ws.Range("O10").Value = 0.1752915379534001
$ws.Range("P10").Value = 0.1752915379534001
$ws.Range("Q10").Value = 2.242340281197333
$ws.Range("R10").Value = 20.181062530776
$ws.Range("S10").Value = 0.0006010982859382968
$ws.Range("T10").Value = 0.0006010982859382968
$ws.Range("I11").Value = 0.003429134645952472
$ws.Range("J11").Value = 0.003429134645952472
$ws.Range("O11").Value = 0.07888758308485012
$ws.Range("P11").Value = 0.07888758308485012
$ws.Range("Q11").Value = 1.009134880683667
$ws.Range("R11").Value = 9.082213926153
$ws.Range("S11").Value = 0.0002705161442917137
$ws.Range("T11").Value = 0.0002705161442917137
$ws.Range("I12").Value = 0.003429134645952472
$ws.Range("J12").Value = 0.003429134645952472
$ws.Range("M12").Value = 2.098187333333334
$ws.Range("N12").Value = 6.294562000000001
$ws.Range("O12").Value = 0.02291182502385553
$ws.Range("P12").Value = 0.02291182502385553
$ws.Range("Q12").Value = 0.2930894940313334
$ws.Range("R12").Value = 2.637805446282
$ws.Range("S12").Value = 0.00007856773299130381
$ws.Range("T12").Value = 0.00007856773299130381
$ws.Range("I13").Value = 0.003429134645952472
$ws.Range("J13").Value = 0.003429134645952472
$ws.Range("M13").Value = 66.20156266666667
$ws.Range("N13").Value = 198.604688
$ws.Range("O13").Value = 0.7229090539378943
$ws.Range("P13").Value = 0.7229090539378942
$ws.Range("Q13").Value = 9.247497684218668
$ws.Range("R13").Value = 83.227479157968
$ws.Range("S13").Value = 0.002478952482731157
$ws.Range("T13").Value = 0.002478952482731157
$ws.Range("G14").Value = 37.47815466666666
$ws.Range("H14").Value = 112.434464
$ws.Range("I14").Value = 0.9200400798487472
$ws.Range("J14").Value = 0.9200400798487472
$ws.Range("M14").Value = 16.05260533333333
$ws.Range("N14").Value = 48.157816
$ws.Range("O14").Value = 0.1752915379534001
$ws.Range("P14").Value = 0.1752915379534001
$ws.Range("Q14").Value = 601.6220254856248
$ws.Range("R14").Value = 5414.598229370623
$ws.Range("S14").Value = 0.1612752405754559
$ws.Range("T14").Value = 0.1612752405754559
$ws.Range("G15").Value = 37.47815466666666
$ws.Range("H15").Value = 112.434464
$ws.Range("I15").Value = 0.9200400798487472
$ws.Range("J15").Value = 0.9200400798487472
$ws.Range("O15").Value = 0.07888758308485012
$ws.Range("P15").Value = 0.07888758308485012
$ws.Range("Q15").Value = 270.7518461831857
$ws.Range("R15").Value = 2436.766615648672
$ws.Range("S15").Value = 0.07257973824046018
$ws.Range("T15").Value = 0.07257973824046018
$ws.Range("G16").Value = 37.47815466666666
$ws.Range("H16").Value = 112.434464
$ws.Range("I16").Value = 0.9200400798487472
$ws.Range("J16").Value = 0.9200400798487472
$ws.Range("M16").Value = 2.098187333333334
$ws.Range("N16").Value = 6.294562000000001
$ws.Range("O16").Value = 0.02291182502385553
$ws.Range("P16").Value = 0.02291182502385553
$ws.Range("Q16").Value = 78.63618939830755
$ws.Range("R16").Value = 707.7257045847681
$ws.Range("S16").Value = 0.02107979732442856
$ws.Range("T16").Value = 0.02107979732442856
$ws.Range("G17").Value = 37.47815466666666
$ws.Range("H17").Value = 112.434464
$ws.Range("I17").Value = 0.9200400798487472
$ws.Range("J17").Value = 0.9200400798487472
$ws.Range("M17").Value = 66.20156266666667
$ws.Range("N17").Value = 198.604688
$ws.Range("O17").Value = 0.7229090539378943
$ws.Range("P17").Value = 0.7229090539378942
$ws.Range("Q17").Value = 2481.112404796359
$ws.Range("R17").Value = 22330.01164316723
$ws.Range("S17").Value = 0.6651053037084026
$ws.Range("T17").Value = 0.6651053037084025
